# Applies the "Updated DJ app and added SAA video" edit:
#   1) "... Saturday, September 22nd, 2018" (Application Process line)
#        -> "... Saturday, September 29nd, 2018"   (2nd digit 2 -> 9)
#   2) "All applications/forms are to be RECEIVED by: Saturday, September
#      22nd, 2018 @ 11:59 pm." (big red headline)
#        -> "... Saturday, September 29nd, 2018 @ 11:59 pm."  (2nd digit 2 -> 9)
#      and the (previously misplaced, in an empty table cell) "_GoBack"
#      bookmark is moved to sit right after the edited digit run, i.e.
#      immediately before the superscript "nd" run.

$d = $word.ActiveDocument

function FindRange($searchText, $startPos, $endPos) {
    $r = $d.Content.Duplicate
    $r.Start = $startPos
    $r.End = $endPos
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) { return $r }
    return $null
}

# ---- Edit 1: "Application Process: ... Saturday, September 22nd, 2018" ----
$anchor1 = FindRange "Application Process:" 0 $d.Content.End
$target1 = FindRange "22nd" $anchor1.End ($anchor1.End + 300)

# Only the first "2" of "22" (the second character overall) changes to "9";
# the rest ("2nd") is untouched.
$digit1 = $d.Range($target1.Start + 1, $target1.Start + 2)
$digit1.Text = "9"

# ---- Edit 2: "... RECEIVED by: Saturday, September 22nd, 2018 @ 11:59 pm." ----
$anchor2 = FindRange "RECEIVED" 0 $d.Content.End
$target2 = FindRange "22nd" $anchor2.End ($anchor2.End + 300)

$digit2 = $d.Range($target2.Start + 1, $target2.Start + 2)
$digit2.Text = "9"

# The edit above shifts the end of the story by 0 (1 char replaced by 1 char),
# so the position right after the edited digit (and before the superscript
# "nd") is still $target2.Start + 2.
$bookmarkPos = $target2.Start + 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)

# Move "_GoBack": drop the old one (empty table cell earlier in the doc) and
# (re)anchor it at the date that was just edited, matching Word's habit of
# leaving "_GoBack" at the location of the most recent edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
